# Commit: "Disabling TC21 in PROD"
# The sanity master executor sheet no longer runs TC21_Verify_Header
# (Verify all Header content) for PROD, so the whole row describing that
# testcase is removed. All rows below it shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

# Row 21 holds TC21_Verify_Header / "Verify  all Header content".
# Deleting the entire row shifts rows 22:25 up to 21:24.
$ws.Rows.Item(21).Delete()

# Matches the saved selection state in the authored workbook.
$ws.Range("C22").Select()
